# Added ESG CO2 emissions
# - Remove the stale _xlchart.v1.* hidden defined names (leftover chart
#   scaffolding no longer referenced).
# - Add two new columns (H: RIC, I: CO2Ratio2018) with Refinitiv RIC
#   tickers and 2018 CO2 ratio figures, matched in for the airlines that
#   have a match.
# - Update selection to H24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove legacy hidden chart-tracking defined names ---
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}

# --- New column headers ---
$ws.Range("H1").Value = "RIC"
$ws.Range("I1").Value = "CO2Ratio2018"

# --- Row 2: EZY / easyJet ---
$ws.Range("H2").Value = "EZJ.L"
$ws.Range("I2").Value = 989.08

# --- Row 3: RYR / Ryanair ---
$ws.Range("H3").Value = "RYA.I"
$ws.Range("I3").Value = 1221.9000000000001

# --- Row 4: EXS / Jet2.com ---
$ws.Range("H4").Value = "DTG.L"

# --- Row 5: DLH / Lufthansa ---
$ws.Range("H5").Value = "LHAG.DE"
$ws.Range("I5").Value = 847.22

# --- Row 6: BEE / Flybe ---
$ws.Range("H6").Value = "FLYB.L^C19"

# --- Row 7: TOM / Thomsonfly ---
$ws.Range("H7").Value = "TUIGn.DE"
$ws.Range("I7").Value = 369.42

# --- Row 8: SHT / British Airways Shuttle ---
$ws.Range("H8").Value = "BAY.L^A11"

# --- Row 9: KLM / KLM Royal Dutch Airlines ---
$ws.Range("H9").Value = "AIRF.PA"
$ws.Range("I9").Value = 920.51

# --- Row 10: SAS / Scandinavian Airlines System ---
$ws.Range("H10").Value = "SAS.ST"
$ws.Range("I10").Value = 947.62

# --- Row 11: EIN / Aer Lingus ---
$ws.Range("H11").Value = "ICAG.L"
$ws.Range("I11").Value = 1074.49

# --- Row 14: QTR / Qatar Airways ---
$ws.Range("H14").Value = "AIRF.PA"
$ws.Range("I14").Value = 920.51

# --- Row 18: SWR / Swiss International Airlines ---
$ws.Range("H18").Value = "LHAG.DE"
$ws.Range("I18").Value = 847.22

# --- Row 23: NAX / Norwegian Air Shuttle ---
$ws.Range("H23").Value = "FIA1S.HE"
$ws.Range("I23").Value = 1004.68

# --- Update active selection ---
[void]$ws.Range("H24").Select()
